$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha/Volumen/Precio mínimo/Precio máximo/Precio promedio ponderado/
# Origen/Precio $/Kg values between row 2 <-> row 5 and row 6 <-> row 7.

function Swap-Rows($ws, $r1, $r2) {
    $cols = @("D", "J", "K", "L", "M", "O", "P")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

Swap-Rows $ws 2 5
Swap-Rows $ws 6 7
